# Refresh Leve profit figures (currentAveragePrice* / LevePrice* / LeveProfit*)
# pulled by the scheduled Sheets runner. Values only -- no formulas in this workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70: Consecrating Congregation / Holy Water
$ws.Range("H70").Value = 2425.8572
$ws.Range("I70").Value = 1989.4
$ws.Range("J70").Value = 2668.3333
$ws.Range("K70").Value = 5968.200000000001
$ws.Range("L70").Value = 8004.999899999999
$ws.Range("M70").Value = -5698.200000000001
$ws.Range("N70").Value = -8544.999899999999
# Row 73: Curbing the Contagion (L) / Holy Water
$ws.Range("H73").Value = 2425.8572
$ws.Range("I73").Value = 1989.4
$ws.Range("J73").Value = 2668.3333
$ws.Range("K73").Value = 5968.200000000001
$ws.Range("L73").Value = 8004.999899999999
$ws.Range("M73").Value = -5032.200000000001
$ws.Range("N73").Value = -9876.999899999999
# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 2832.7273
$ws.Range("J112").Value = 2989.9
$ws.Range("L112").Value = 8969.700000000001
$ws.Range("N112").Value = -11185.7
# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 4737.7046
$ws.Range("I132").Value = 3085.4443
$ws.Range("J132").Value = 12172.875
$ws.Range("K132").Value = 9256.332900000001
$ws.Range("L132").Value = 36518.625
$ws.Range("M132").Value = -6726.332900000001
$ws.Range("N132").Value = -41578.625
# Row 135: For Tired Minds / Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 889.69446
$ws.Range("I135").Value = 897.13794
$ws.Range("J135").Value = 858.8570999999999
$ws.Range("K135").Value = 8074.241459999999
$ws.Range("L135").Value = 7729.7139
$ws.Range("M135").Value = -5539.241459999999
$ws.Range("N135").Value = -12799.7139
# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2019.0392
$ws.Range("I138").Value = 1078.625
$ws.Range("J138").Value = 3602.8948
$ws.Range("K138").Value = 3235.875
$ws.Range("L138").Value = 10808.6844
$ws.Range("M138").Value = 1904.125
$ws.Range("N138").Value = -21088.6844
# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 1640.1555
$ws.Range("I141").Value = 1662.091
$ws.Range("K141").Value = 4986.272999999999
$ws.Range("M141").Value = 193.7270000000008

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 15120.023
$ws.Range("I32").Value = 14584.734
$ws.Range("K32").Value = 14584.734
$ws.Range("M32").Value = -14297.734
# Row 43: They've Got Legs / Steel Sabatons
$ws.Range("H43").Value = 19990
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 19990
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 19990
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -20616
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 1762.1578
$ws.Range("I61").Value = 1252.4286
$ws.Range("K61").Value = 1252.4286
$ws.Range("M61").Value = -1040.4286
# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 29890.432
$ws.Range("I74").Value = 32203.617
$ws.Range("J74").Value = 3674.3333
$ws.Range("K74").Value = 32203.617
$ws.Range("L74").Value = 3674.3333
$ws.Range("M74").Value = -31329.617
$ws.Range("N74").Value = -5422.3333
# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 29890.432
$ws.Range("I77").Value = 32203.617
$ws.Range("J77").Value = 3674.3333
$ws.Range("K77").Value = 161018.085
$ws.Range("L77").Value = 18371.6665
$ws.Range("M77").Value = -156650.085
$ws.Range("N77").Value = -27107.6665
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 22431.896
$ws.Range("I132").Value = 24186.955
$ws.Range("J132").Value = 3126.25
$ws.Range("K132").Value = 72560.86500000001
$ws.Range("L132").Value = 9378.75
$ws.Range("M132").Value = -70030.86500000001
$ws.Range("N132").Value = -14438.75
# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 1762.1578
$ws.Range("I136").Value = 1252.4286
$ws.Range("K136").Value = 3757.2858
$ws.Range("M136").Value = -1207.2858

$ws = $wb.Worksheets.Item("BSM")
# Row 94: High Steal / High Steel Nugget
$ws.Range("H94").Value = 1706.08
$ws.Range("I94").Value = 1287.6666
$ws.Range("J94").Value = 2782
$ws.Range("K94").Value = 1287.6666
$ws.Range("L94").Value = 2782
$ws.Range("M94").Value = -836.6666
$ws.Range("N94").Value = -3684
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 2815.1516
$ws.Range("I134").Value = 2521.1035
$ws.Range("J134").Value = 4947
$ws.Range("K134").Value = 7563.310500000001
$ws.Range("L134").Value = 14841
$ws.Range("M134").Value = -5028.310500000001
$ws.Range("N134").Value = -19911

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 3426.2144
$ws.Range("I31").Value = 3053.6365
$ws.Range("J31").Value = 4792.3335
$ws.Range("K31").Value = 3053.6365
$ws.Range("L31").Value = 4792.3335
$ws.Range("M31").Value = -2758.6365
$ws.Range("N31").Value = -5382.3335
# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 3426.2144
$ws.Range("I34").Value = 3053.6365
$ws.Range("J34").Value = 4792.3335
$ws.Range("K34").Value = 3053.6365
$ws.Range("L34").Value = 4792.3335
$ws.Range("M34").Value = -2851.6365
$ws.Range("N34").Value = -5196.3335
# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 34557.355
$ws.Range("I58").Value = 42259.6
$ws.Range("J58").Value = 2464.6667
$ws.Range("K58").Value = 42259.6
$ws.Range("L58").Value = 2464.6667
$ws.Range("M58").Value = -42056.6
$ws.Range("N58").Value = -2870.6667
# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 34557.355
$ws.Range("I136").Value = 42259.6
$ws.Range("J136").Value = 2464.6667
$ws.Range("K136").Value = 126778.8
$ws.Range("L136").Value = 7394.000100000001
$ws.Range("M136").Value = -124228.8
$ws.Range("N136").Value = -12494.0001

$ws = $wb.Worksheets.Item("CUL")
# Row 92: Oh No Udon / Gyr Abanian Flour
$ws.Range("H92").Value = 229.76596
$ws.Range("I92").Value = 273.83334
$ws.Range("K92").Value = 821.5000200000001
$ws.Range("M92").Value = 426.4999799999999

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 2570.394
$ws.Range("I80").Value = 2005.4667
$ws.Range("J80").Value = 3041.1667
$ws.Range("K80").Value = 2005.4667
$ws.Range("L80").Value = 3041.1667
$ws.Range("M80").Value = -1007.4667
$ws.Range("N80").Value = -5037.1667
# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 2570.394
$ws.Range("I83").Value = 2005.4667
$ws.Range("J83").Value = 3041.1667
$ws.Range("K83").Value = 10027.3335
$ws.Range("L83").Value = 15205.8335
$ws.Range("M83").Value = -5035.333500000001
$ws.Range("N83").Value = -25189.8335
# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 1816.0322
$ws.Range("I122").Value = 1768.862
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 5306.586
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -2856.586
$ws.Range("N122").Value = -12400
# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 4697.7827
$ws.Range("I126").Value = 4558.6113
$ws.Range("J126").Value = 5198.8
$ws.Range("K126").Value = 13675.8339
$ws.Range("L126").Value = 15596.4
$ws.Range("M126").Value = -11205.8339
$ws.Range("N126").Value = -20536.4
# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 35056.098
$ws.Range("I132").Value = 42643.76
$ws.Range("J132").Value = 3440.8333
$ws.Range("K132").Value = 127931.28
$ws.Range("L132").Value = 10322.4999
$ws.Range("M132").Value = -125401.28
$ws.Range("N132").Value = -15382.4999

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore / Hard Leather
$ws.Range("H16").Value = 425
$ws.Range("I16").Value = 425
$ws.Range("K16").Value = 425
$ws.Range("M16").Value = -255
# Row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 3891
$ws.Range("I68").Value = 3676.889
$ws.Range("K68").Value = 3676.889
$ws.Range("M68").Value = -2927.889
# Row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 3891
$ws.Range("I71").Value = 3676.889
$ws.Range("K71").Value = 18384.445
$ws.Range("M71").Value = -14640.445
# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 1887.8636
$ws.Range("I93").Value = 1577.0714
$ws.Range("J93").Value = 2431.75
$ws.Range("K93").Value = 1577.0714
$ws.Range("L93").Value = 2431.75
$ws.Range("M93").Value = -329.0714
$ws.Range("N93").Value = -4927.75
# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 213445.83
$ws.Range("I122").Value = 5692.8276
$ws.Range("J122").Value = 530542.5
$ws.Range("K122").Value = 17078.4828
$ws.Range("L122").Value = 1591627.5
$ws.Range("M122").Value = -14628.4828
$ws.Range("N122").Value = -1596527.5
# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 28298.688
$ws.Range("I132").Value = 32748.775
$ws.Range("J132").Value = 6048.25
$ws.Range("K132").Value = 98246.32500000001
$ws.Range("L132").Value = 18144.75
$ws.Range("M132").Value = -95716.32500000001
$ws.Range("N132").Value = -23204.75
# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 4528.0835
$ws.Range("I136").Value = 4093.3333
$ws.Range("J136").Value = 5832.3335
$ws.Range("K136").Value = 12279.9999
$ws.Range("L136").Value = 17497.0005
$ws.Range("M136").Value = -9729.999899999999
$ws.Range("N136").Value = -22597.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 62: Pride Up in Smoke / Rainbow Cloth
$ws.Range("H62").Value = 213395.6
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 213395.6
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 213395.6
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -214643.6
# Row 65: Desperate for Diversionaries (L) / Rainbow Cloth
$ws.Range("H65").Value = 213395.6
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 213395.6
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 1066978
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -1073218
# Row 81: Where the Dragonflies, the Net Catches / Crawler Silk
$ws.Range("H81").Value = 6015.6665
$ws.Range("I81").Value = 1356.5
$ws.Range("J81").Value = 15334
$ws.Range("K81").Value = 2713
$ws.Range("L81").Value = 30668
$ws.Range("M81").Value = -1652
$ws.Range("N81").Value = -32790
# Row 84: To Kill a Dragon on Nameday (L) / Crawler Silk
$ws.Range("H84").Value = 6015.6665
$ws.Range("I84").Value = 1356.5
$ws.Range("J84").Value = 15334
$ws.Range("K84").Value = 13565
$ws.Range("L84").Value = 153340
$ws.Range("M84").Value = -8261
$ws.Range("N84").Value = -163948
# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 29663.695
$ws.Range("I122").Value = 1835.5
$ws.Range("J122").Value = 127062.375
$ws.Range("K122").Value = 5506.5
$ws.Range("L122").Value = 381187.125
$ws.Range("M122").Value = -3056.5
$ws.Range("N122").Value = -386087.125
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 1436.3334
$ws.Range("I136").Value = 1313.8572
$ws.Range("K136").Value = 3941.5716
$ws.Range("M136").Value = -1391.5716
